$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the data value in D1
$ws.Range("D1").Value = 511207946

# Update the active selection to D1 (single cell)
$ws.Range("D1").Select()
